$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 127/128, shifting the existing rows 127-151 down to 129-153
$ws.Range("A127:A128").EntireRow.Insert()

# New row 127: Espárragos, Sin especificar, Primera, Provincia de Linares
$ws.Cells.Item(127, 1).Value  = 9
$ws.Cells.Item(127, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(127, 3).Value  = "Metropolitana"
$ws.Cells.Item(127, 4).Value  = 44841
$ws.Cells.Item(127, 5).Value  = 13
$ws.Cells.Item(127, 6).Value  = 300000000
$ws.Cells.Item(127, 7).Value  = "Espárragos"
$ws.Cells.Item(127, 8).Value  = "Sin especificar"
$ws.Cells.Item(127, 9).Value  = "Primera"
$ws.Cells.Item(127, 10).Value = 540
$ws.Cells.Item(127, 11).Value = 1400
$ws.Cells.Item(127, 12).Value = 1400
$ws.Cells.Item(127, 13).Value = 1400
$ws.Cells.Item(127, 14).Value = "`$/kilo"
$ws.Cells.Item(127, 15).Value = "Provincia de Linares"
$ws.Cells.Item(127, 16).Value = 1400
$ws.Cells.Item(127, 17).Value = 1
$ws.Cells.Item(127, 18).Value = "Hortaliza"

# New row 128: Espárragos, Sin especificar, Segunda, Provincia de Linares
$ws.Cells.Item(128, 1).Value  = 9
$ws.Cells.Item(128, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(128, 3).Value  = "Metropolitana"
$ws.Cells.Item(128, 4).Value  = 44841
$ws.Cells.Item(128, 5).Value  = 13
$ws.Cells.Item(128, 6).Value  = 300000000
$ws.Cells.Item(128, 7).Value  = "Espárragos"
$ws.Cells.Item(128, 8).Value  = "Sin especificar"
$ws.Cells.Item(128, 9).Value  = "Segunda"
$ws.Cells.Item(128, 10).Value = 360
$ws.Cells.Item(128, 11).Value = 1100
$ws.Cells.Item(128, 12).Value = 1100
$ws.Cells.Item(128, 13).Value = 1100
$ws.Cells.Item(128, 14).Value = "`$/kilo"
$ws.Cells.Item(128, 15).Value = "Provincia de Linares"
$ws.Cells.Item(128, 16).Value = 1100
$ws.Cells.Item(128, 17).Value = 1
$ws.Cells.Item(128, 18).Value = "Hortaliza"
